$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_6_3_0", 0.303141395600051, 0.4988412828114286, 0.4995191805889216, 0.5018448103859863, 0.7712171673774719, 0.3687612116336823, 0.6017147898674011, 0.4783864915370941),
    @("model_6_3_7", 0.3159747985577135, 0.3607374873541642, 0.429186821329636, 0.4043864014648071, 0.7570143342018127, 0.4703803658485413, 0.686273455619812, 0.5719773769378662),
    @("model_6_3_10", 0.3169217561628715, 0.3546669916352498, 0.4317826564056003, 0.4034534640511704, 0.7559663057327271, 0.4748471677303314, 0.6831525564193726, 0.5728732943534851),
    @("model_6_3_8", 0.3171422555387216, 0.3570529563849864, 0.4360443296713432, 0.4069323374104857, 0.7557222247123718, 0.4730915427207947, 0.6780288219451904, 0.5695325136184692),
    @("model_6_3_9", 0.3173245238764135, 0.3560550958577438, 0.4355999165622462, 0.4062654451725211, 0.7555206418037415, 0.4738257825374603, 0.6785631775856018, 0.5701728463172913),
    @("model_6_3_11", 0.3180445781540041, 0.3543615729918129, 0.438363463906168, 0.4072065677096386, 0.7547237277030945, 0.475071907043457, 0.6752406358718872, 0.5692691802978516),
    @("model_6_3_17", 0.318643678875937, 0.3525794886581808, 0.4380591566365143, 0.4063044414755439, 0.7540607452392578, 0.4763831496238708, 0.6756064891815186, 0.5701355338096619),
    @("model_6_3_12", 0.3190953711377051, 0.3531386373958869, 0.4448346445346285, 0.4105232862844839, 0.7535607218742371, 0.4759716987609863, 0.667460560798645, 0.5660840272903442),
    @("model_6_3_5", 0.3191670294576114, 0.3832991017680522, 0.4524027780801865, 0.4272164325641168, 0.7534814476966858, 0.4537791609764099, 0.658361554145813, 0.5500532984733582),
    @("model_6_3_19", 0.3194940962686101, 0.3520374895338033, 0.442724380841599, 0.4088332120927317, 0.7531195282936096, 0.476781964302063, 0.6699976325035095, 0.5677070617675781),
    @("model_6_3_13", 0.3198536368724223, 0.3525496042924252, 0.4496448435450476, 0.4131184489468089, 0.7527216672897339, 0.4764052033424377, 0.661677360534668, 0.5635918974876404),
    @("model_6_3_20", 0.3200289011371344, 0.3517003242606148, 0.4463018429756695, 0.4108041101333579, 0.75252765417099, 0.4770301282405853, 0.6656965017318726, 0.5658143758773804),
    @("model_6_3_18", 0.3200347206804754, 0.3523932313192837, 0.4464253886762505, 0.4111581653860413, 0.7525211572647095, 0.4765202701091766, 0.6655480265617371, 0.5654743313789368),
    @("model_6_3_21", 0.320287891242039, 0.3516711970285089, 0.4479716184395525, 0.4117764771821546, 0.7522410154342651, 0.4770515263080597, 0.6636890172958374, 0.5648806095123291),
    @("model_6_3_23", 0.3204000118116243, 0.3518414532697464, 0.4480947801150216, 0.4119180348423229, 0.7521169781684875, 0.4769262075424194, 0.6635409593582153, 0.5647446513175964),
    @("model_6_3_1", 0.3204141950533059, 0.4916886253844577, 0.5254222704883971, 0.5142040027265902, 0.7521012425422668, 0.3740242719650269, 0.5705721378326416, 0.4665177464485168),
    @("model_6_3_24", 0.3205897271361512, 0.3517187819379062, 0.4493301945159272, 0.4125959367908754, 0.7519069314002991, 0.4770165085792542, 0.6620556116104126, 0.5640936493873596),
    @("model_6_3_22", 0.3207374337245791, 0.3516308571867313, 0.450480200221421, 0.4132377625290103, 0.751743495464325, 0.4770811796188354, 0.6606730222702026, 0.5634772777557373),
    @("model_6_3_14", 0.3207686169936688, 0.3519339938114173, 0.4553848918737579, 0.4162503400128897, 0.751708984375, 0.4768581390380859, 0.6547762155532837, 0.5605842471122742),
    @("model_6_3_16", 0.3216455401452794, 0.3526329130773129, 0.4580643158347577, 0.418112668534668, 0.7507385015487671, 0.4763438999652863, 0.6515548229217529, 0.5587958097457886),
    @("model_6_3_15", 0.3216700479940987, 0.3536588227386407, 0.459233549520738, 0.4192173674293385, 0.7507113814353943, 0.4755889475345612, 0.6501491069793701, 0.5577349662780762),
    @("model_6_3_6", 0.323365097695361, 0.3702337151015238, 0.4816711522452284, 0.4391602177090889, 0.7488355040550232, 0.4633928835391998, 0.6231729984283447, 0.538583517074585),
    @("model_6_3_2", 0.3253907268338032, 0.4307462015770821, 0.5359567141980421, 0.4956896250273592, 0.7465936541557312, 0.418866753578186, 0.5579068660736084, 0.4842974245548248),
    @("model_6_3_3", 0.3275685635079297, 0.4247068090509418, 0.5370802951286908, 0.4939015578724805, 0.7441834807395935, 0.4233106374740601, 0.5565560460090637, 0.4860145449638367),
    @("model_6_3_4", 0.3276792161835731, 0.4183716800333475, 0.5331016172972414, 0.4889881147519561, 0.7440609931945801, 0.427972137928009, 0.5613394975662231, 0.4907329678535461)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

Write-Output "applied"
